# EventResult_Data.xlsx update — "3 new EventTrigger Update"
#
# 1) Row 30 col B: rename Dangerous_Mission_2_4_B -> Dangerous_Mission_2_4_A
#    (this trigger id is being reused/normalized now that more "_A" rows
#    follow it).
# 2) Append three brand-new EventTrigger rows (407/408/409) that reuse the
#    same B/E/F values as row 30/29, each with its own Chinese-text
#    condition in column C and an incrementing D (10, 11, 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 30: Dangerous_Mission_2_4_B -> Dangerous_Mission_2_4_A
# ---------------------------------------------------------------------
$ws.Range("B30").Value = "Dangerous_Mission_2_4_A"

# ---------------------------------------------------------------------
# 2) New rows 31-33, formatted like row 30 (same row height + per-column
#    alignment/wrap) with 3 new trigger conditions.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=31; A=407; C="角色到了";         D=10 },
    @{ Row=32; A=408; C="建造了一个弹药箱"; D=11 },
    @{ Row=33; A=409; C="角色开打";         D=12 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A" + $row + ":F" + $row).RowHeight = 27.6

    $ws.Range("A" + $row).Value = $r.A
    $ws.Range("A" + $row).HorizontalAlignment = -4108
    $ws.Range("A" + $row).VerticalAlignment = -4108

    $ws.Range("B" + $row).Value = "Dangerous_Mission_2_4_A"
    $ws.Range("B" + $row).HorizontalAlignment = -4108
    $ws.Range("B" + $row).VerticalAlignment = -4108

    $ws.Range("C" + $row).Value = $r.C
    $ws.Range("C" + $row).HorizontalAlignment = -4131
    $ws.Range("C" + $row).VerticalAlignment = -4108
    $ws.Range("C" + $row).WrapText = $true

    $ws.Range("D" + $row).Value = $r.D
    $ws.Range("D" + $row).HorizontalAlignment = -4108
    $ws.Range("D" + $row).VerticalAlignment = -4108

    $ws.Range("E" + $row).Value = "[305]"
    $ws.Range("E" + $row).HorizontalAlignment = -4108
    $ws.Range("E" + $row).VerticalAlignment = -4108

    $ws.Range("F" + $row).Value = "Dangerous_Mission_0_1"
    $ws.Range("F" + $row).HorizontalAlignment = -4108
    $ws.Range("F" + $row).VerticalAlignment = -4108
    $ws.Range("F" + $row).WrapText = $true
}

# ---------------------------------------------------------------------
# View state: scroll down to the new rows and select F29:F33 like the
# author did after adding them.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 55
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F29:F33").Select()
